$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("A2").Value = "FEP0000001"
$ws.Range("C2").Value = "27/01/2025 11:22:45"
$ws.Range("D2").Value = "27/01/2025 11:24:27"
$ws.Range("G2").Value = "00:01:01"

# Row 3 updates
$ws.Range("A3").Value = "FEP0000002"
$ws.Range("C3").Value = "27/01/2025 09:23:12"
$ws.Range("D3").Value = "27/01/2025 11:24:46"
$ws.Range("G3").Value = "02:01:01"
$ws.Range("H3").Value = 10000

# Row 4 new entry
$ws.Range("A4").Value = "FEP0000003"
$ws.Range("B4").Value = "CCC01"
$ws.Range("C4").Value = "26/01/2025 23:26:22"
$ws.Range("D4").Value = "27/01/2025 11:28:02"
$ws.Range("E4").Value = "Moto"
$ws.Range("F4").Value = 7000
$ws.Range("G4").Value = "12:01:01"
$ws.Range("H4").Value = 7750
